# Update countries & provincias Spain
# Refreshes the COVID "Pais" dashboard data from the 17:29 snapshot to the
# 18:46 snapshot: updates the "data as of" timestamp, refreshes the numeric
# columns (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for every country whose figures
# changed, and re-labels the handful of rows whose country swapped rank
# with its neighbour (Singapur/Israel, Corea del Sur/Chequia,
# Madagascar/Guayana Francesa, Sri Lanka/Libano,
# Sahara Occidental/Bonaire-San Eustaquio y Saba).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 18:46"

# Row 4
$ws.Cells.Item(4, 2).Value = 3721626
$ws.Cells.Item(4, 3).Value = 26601
$ws.Cells.Item(4, 4).Value = 1682009
$ws.Cells.Item(4, 5).Value = 1898161
$ws.Cells.Item(4, 7).Value = 338
$ws.Cells.Item(4, 8).Value = 141456

# Row 5
$ws.Cells.Item(5, 2).Value = 2021834
$ws.Cells.Item(5, 3).Value = 7096
$ws.Cells.Item(5, 5).Value = 578062
$ws.Cells.Item(5, 7).Value = 175
$ws.Cells.Item(5, 8).Value = 76997

# Row 6
$ws.Cells.Item(6, 2).Value = 1037249
$ws.Cells.Item(6, 3).Value = 31612
$ws.Cells.Item(6, 4).Value = 652582
$ws.Cells.Item(6, 5).Value = 358394
$ws.Cells.Item(6, 7).Value = 664
$ws.Cells.Item(6, 8).Value = 26273

# Row 12
$ws.Cells.Item(12, 2).Value = 307335
$ws.Cells.Item(12, 3).Value = 1400
$ws.Cells.Item(12, 7).Value = 4
$ws.Cells.Item(12, 8).Value = 28420

# Row 17
$ws.Cells.Item(17, 2).Value = 243967
$ws.Cells.Item(17, 3).Value = 231
$ws.Cells.Item(17, 4).Value = 196483
$ws.Cells.Item(17, 5).Value = 12456
$ws.Cells.Item(17, 7).Value = 11
$ws.Cells.Item(17, 8).Value = 35028

# Row 24
$ws.Cells.Item(24, 2).Value = 109516
$ws.Cells.Item(24, 3).Value = 252
$ws.Cells.Item(24, 4).Value = 96623
$ws.Cells.Item(24, 5).Value = 4058
$ws.Cells.Item(24, 7).Value = 8
$ws.Cells.Item(24, 8).Value = 8835

# Row 45
$ws.Cells.Item(45, 1).Value = "Israel"
$ws.Cells.Item(45, 2).Value = 47459
$ws.Cells.Item(45, 3).Value = 1400
$ws.Cells.Item(45, 4).Value = 20744
$ws.Cells.Item(45, 5).Value = 26323
$ws.Cells.Item(45, 7).Value = 8
$ws.Cells.Item(45, 8).Value = 392

# Row 46
$ws.Cells.Item(46, 1).Value = "Singapur"
$ws.Cells.Item(46, 2).Value = 47453
$ws.Cells.Item(46, 3).Value = 327
$ws.Cells.Item(46, 4).Value = 43577
$ws.Cells.Item(46, 5).Value = 3849
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 27

# Row 50
$ws.Cells.Item(50, 5).Value = 4151
$ws.Cells.Item(50, 7).Value = 3
$ws.Cells.Item(50, 8).Value = 124

# Row 60
$ws.Cells.Item(60, 2).Value = 21948
$ws.Cells.Item(60, 3).Value = 593
$ws.Cells.Item(60, 4).Value = 15430
$ws.Cells.Item(60, 5).Value = 5461
$ws.Cells.Item(60, 7).Value = 5
$ws.Cells.Item(60, 8).Value = 1057

# Row 68
$ws.Cells.Item(68, 1).Value = "Chequia"
$ws.Cells.Item(68, 2).Value = 13682
$ws.Cells.Item(68, 3).Value = 70
$ws.Cells.Item(68, 4).Value = 8725
$ws.Cells.Item(68, 5).Value = 4599
$ws.Cells.Item(68, 7).Value = 3
$ws.Cells.Item(68, 8).Value = 358

# Row 69
$ws.Cells.Item(69, 1).Value = "Corea del Sur"
$ws.Cells.Item(69, 2).Value = 13672
$ws.Cells.Item(69, 3).Value = 60
$ws.Cells.Item(69, 4).Value = 12460
$ws.Cells.Item(69, 5).Value = 919
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 293

# Row 87
$ws.Cells.Item(87, 5).Value = 6220
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = 52

# Row 91
$ws.Cells.Item(91, 1).Value = "Guayana Francesa"
$ws.Cells.Item(91, 2).Value = 6509
$ws.Cells.Item(91, 3).Value = 116
$ws.Cells.Item(91, 4).Value = 3932
$ws.Cells.Item(91, 5).Value = 2543
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 34

# Row 92
$ws.Cells.Item(92, 1).Value = "Madagascar"
$ws.Cells.Item(92, 2).Value = 6467
$ws.Cells.Item(92, 3).Value = 378
$ws.Cells.Item(92, 4).Value = 3108
$ws.Cells.Item(92, 5).Value = 3305
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = 54

# Row 96
$ws.Cells.Item(96, 2).Value = 5409
$ws.Cells.Item(96, 3).Value = 124
$ws.Cells.Item(96, 4).Value = 4333
$ws.Cells.Item(96, 5).Value = 965

# Row 101
$ws.Cells.Item(101, 2).Value = 3964
$ws.Cells.Item(101, 3).Value = 25
$ws.Cells.Item(101, 5).Value = 2396
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 194

# Row 112
$ws.Cells.Item(112, 1).Value = "Libano"
$ws.Cells.Item(112, 2).Value = 2700
$ws.Cells.Item(112, 3).Value = 101
$ws.Cells.Item(112, 4).Value = 1485
$ws.Cells.Item(112, 5).Value = 1175
$ws.Cells.Item(112, 8).Value = 40

# Row 113
$ws.Cells.Item(113, 1).Value = "Sri Lanka"
$ws.Cells.Item(113, 2).Value = 2689
$ws.Cells.Item(113, 3).Value = 2
$ws.Cells.Item(113, 4).Value = 2012
$ws.Cells.Item(113, 5).Value = 666
$ws.Cells.Item(113, 8).Value = 11

# Row 127
$ws.Cells.Item(127, 2).Value = 1688
$ws.Cells.Item(127, 3).Value = 10
$ws.Cells.Item(127, 4).Value = 1219
$ws.Cells.Item(127, 5).Value = 404
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = 65

# Row 137
$ws.Cells.Item(137, 2).Value = 1209
$ws.Cells.Item(137, 3).Value = 3
$ws.Cells.Item(137, 4).Value = 1021
$ws.Cells.Item(137, 5).Value = 178

# Row 148
$ws.Cells.Item(148, 2).Value = 887
$ws.Cells.Item(148, 3).Value = 1
$ws.Cells.Item(148, 4).Value = 800

# Row 180
$ws.Cells.Item(180, 2).Value = 136
$ws.Cells.Item(180, 3).Value = 3
$ws.Cells.Item(180, 5).Value = 4

# Row 214
$ws.Cells.Item(214, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(214, 3).Value = 1
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = 3
$ws.Cells.Item(214, 8).Value = 0

# Row 215
$ws.Cells.Item(215, 1).Value = "Sahara Occidental"
$ws.Cells.Item(215, 2).Value = 10
$ws.Cells.Item(215, 4).Value = 8
$ws.Cells.Item(215, 5).Value = 1
$ws.Cells.Item(215, 8).Value = 1
